# Add Haroun's details as a new row under the existing header row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "عبدالحميد عادل محمود شرادة"
$ws.Range("B2").Value = "harounwaka125@gmail.com"
$ws.Range("C2").Value = "https://github.com/harounwaka125/OOS-Project"

# Turn the email and repo link into real hyperlinks (this also applies the
# built-in "Hyperlink" cell style - underline + theme color)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:harounwaka125@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/harounwaka125/OOS-Project")

# Resize the columns so the new, longer content fits (best-fit widths)
$ws.Columns.Item(1).ColumnWidth = 19.0
$ws.Columns.Item(2).ColumnWidth = 23.166666666666668
$ws.Columns.Item(3).ColumnWidth = 40.0

$ws.Range("C7").Select()
